# Updates the cryptos price/volume table on Sheet1, reflecting the
# latest coinranking.com scrape (includes two coin row swaps), per commit:
# "Updated cryptos list on Mon Feb 12 02:47:29 UTC 2024 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "48.094.75"
$ws.Range("E2").Value = "  +0.74%  "
$ws.Range("D3").Value = "2.502.39"
$ws.Range("E3").Value = "  -0.23%  "
$c = $ws.Range("D4")
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "
$c = $ws.Range("D5")
$c.Value = "'318.73"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.26%  "
$c = $ws.Range("D6")
$c.Value = "'105.80"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -2.18%  "
$c = $ws.Range("D7")
$c.Value = "'0.523"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.12%  "
$c = $ws.Range("D8")
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  -3.98%  "
$c = $ws.Range("D10")
$c.Value = "'38.90"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -3.41%  "
$c = $ws.Range("D11")
$c.Value = "'19.87"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +0.66%  "
$c = $ws.Range("D12")
$c.Value = "'0.0803"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -1.29%  "
$ws.Range("E13").Value = "  -0.83%  "
$c = $ws.Range("D14")
$c.Value = "'7.05"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -1.80%  "
$ws.Range("D15").Value = "2.897.34"
$ws.Range("E15").Value = "  -0.09%  "
$ws.Range("D16").Value = "2.507.63"
$ws.Range("E16").Value = "  -0.09%  "
$c = $ws.Range("D17")
$c.Value = "'0.830"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -2.45%  "
$ws.Range("D18").Value = "47.987.69"
$ws.Range("E18").Value = "  +0.71%  "
$c = $ws.Range("D19")
$c.Value = "'12.95"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -2.88%  "
$c = $ws.Range("D20")
$c.Value = "'2.95"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +6.92%  "
$c = $ws.Range("D21")
$c.Value = "'6.63"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.30%  "
$ws.Range("D22").Value = "0.0₃0934"
$ws.Range("E22").Value = "  -0.80%  "
$c = $ws.Range("D23")
$c.Value = "'71.01"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.06%  "
$c = $ws.Range("D24")
$c.Value = "'271.69"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +9.79%  "
$ws.Range("E25").Value = "  -2.03%  "
$ws.Range("E26").Value = "  -0.05%  "
$c = $ws.Range("D27")
$c.Value = "'25.79"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$c = $ws.Range("D28")
$c.Value = "'2.25"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +2.48%  "
$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c = $ws.Range("D29")
$c.Value = "'0.144"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +1.25%  "
$c = $ws.Range("D30")
$c.Value = "'9.70"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -5.16%  "
$c = $ws.Range("D31")
$c.Value = "'34.56"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -0.82%  "
$c = $ws.Range("D32")
$c.Value = "'49.31"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -1.08%  "
$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$c = $ws.Range("D33")
$c.Value = "'1.00"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("B34").Value = "Celestia"
$ws.Range("C34").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$c = $ws.Range("D34")
$c.Value = "'19.10"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -4.80%  "
$c = $ws.Range("D35")
$c.Value = "'5.28"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -1.81%  "
$c = $ws.Range("D36")
$c.Value = "'0.0776"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -1.08%  "
$c = $ws.Range("D37")
$c.Value = "'1.94"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -1.41%  "
$c = $ws.Range("D38")
$c.Value = "'4.59"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -2.06%  "
$c = $ws.Range("D39")
$c.Value = "'2.86"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -3.83%  "
$ws.Range("E40").Value = "  -1.40%  "
$ws.Range("E41").Value = "  +1.00%  "
$c = $ws.Range("D42")
$c.Value = "'119.84"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +1.10%  "
$c = $ws.Range("D43")
$c.Value = "'21.74"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -2.72%  "
$c = $ws.Range("D44")
$c.Value = "'0.0303"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +1.80%  "
$ws.Range("D45").Value = "2.002.68"
$ws.Range("E45").Value = "  +0.04%  "
$c = $ws.Range("D46")
$c.Value = "'3.20"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +3.33%  "
$c = $ws.Range("D47")
$c.Value = "'1.87"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +4.08%  "
$c = $ws.Range("D49")
$c.Value = "'8.93"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -1.78%  "
$c = $ws.Range("D50")
$c.Value = "'5.18"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.29%  "
$c = $ws.Range("D51")
$c.Value = "'78.75"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +1.71%  "

